$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jun")
$ws.Activate()

# New status rows (87-94), matching the existing table layout (columns A-F).
# Column A already has the running index filled in; copy formatting from the
# row above (row 86) so the new cells pick up the same styles (incl. the
# date number format in column F), then fill in the values.
$rows = @(
  @{ Row = 87; B = "브루트 포스"; C = 2; D = 15665; E = "N과 M(11)";   F = 44217 },
  @{ Row = 88; B = "브루트 포스"; C = 1; D = 15666; E = "N과 M(12)";   F = 44217 },
  @{ Row = 89; B = "브루트 포스"; C = 2; D = 10972; E = "다음 수열";   F = 44217 },
  @{ Row = 90; B = "브루트 포스"; C = 1; D = 10973; E = "이전 수열";   F = 44217 },
  @{ Row = 91; B = "브루트 포스"; C = 1; D = 10974; E = "모든 수열";   F = 44217 },
  @{ Row = 92; B = "브루트 포스"; C = 1; D = 10819; E = "차이를 최대로"; F = 44217 },
  @{ Row = 93; B = "브루트 포스"; C = 1; D = 10971; E = "외판원 순회2"; F = 44217 },
  @{ Row = 94; B = "브루트 포스"; C = 2; D = 6603;  E = "로또";        F = 44217 }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Range("A86:F86").Copy()
  $ws.Range("A" + $row + ":F" + $row).PasteSpecial(-4122)

  $ws.Range("B" + $row).Value = $r.B
  $ws.Range("C" + $row).Value = $r.C
  $ws.Range("D" + $row).Value = $r.D
  $ws.Range("E" + $row).Value = $r.E
  $ws.Range("F" + $row).Value = $r.F
}

$excel.CutCopyMode = 0

# Update the sheet's scroll position / selection to match the latest edit.
$ws.Range("B87:B94").Select()
